# Apply the cryptos-list price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay a TEXT cell even when it looks like a
# number (e.g. "24.47"), without leaving the cells style/NumberFormat altered.
# We flip the cell to Text format just long enough to enter the value, then put
# the original Style back so the saved XML keeps the same style index as before.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

# Row 2
$ws.Range("D2").Value = '29.363.48'
$ws.Range("E2").Value = '  -0.34%  '

# Row 3
$ws.Range("D3").Value = '1.845.99'
$ws.Range("E3").Value = '  -0.19%  '

# Row 4
Set-TextValue "D4" '0.9988'
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
Set-TextValue "D5" '240.08'
$ws.Range("E5").Value = '  -0.73%  '

# Row 6
Set-TextValue "D6" '0.6307'
$ws.Range("E6").Value = '  +0.64%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
Set-TextValue "D8" '0.07538'
$ws.Range("E8").Value = '  +0.14%  '

# Row 9
Set-TextValue "D9" '0.2958'
$ws.Range("E9").Value = '  -0.57%  '

# Row 10
Set-TextValue "D10" '24.47'
$ws.Range("E10").Value = '  +0.64%  '

# Row 11
Set-TextValue "D11" '0.07720'
$ws.Range("E11").Value = '  +0.18%  '

# Row 12
$ws.Range("D12").Value = '1.851.94'
$ws.Range("E12").Value = '  -5.19%  '

# Row 13
$ws.Range("E13").Value = '  -0.09%  '

# Row 14
Set-TextValue "D14" '0.6837'
$ws.Range("E14").Value = '  -0.12%  '

# Row 15
Set-TextValue "D15" '0.000009998'
$ws.Range("E15").Value = '  +2.67%  '

# Row 16
Set-TextValue "D16" '82.83'
$ws.Range("E16").Value = '  -1.10%  '

# Row 17
Set-TextValue "D17" '6.141'
$ws.Range("E17").Value = '  -1.28%  '

# Row 18
$ws.Range("D18").Value = '29.391.10'
$ws.Range("E18").Value = '  -0.90%  '

# Row 19
Set-TextValue "D19" '227.70'
$ws.Range("E19").Value = '  -2.58%  '

# Row 21
Set-TextValue "D21" '0.9998'
$ws.Range("E21").Value = '  +0.01%  '

# Row 22
Set-TextValue "D22" '7.539'
$ws.Range("E22").Value = '  -0.62%  '

# Row 23
$ws.Range("E23").Value = '  +0.04%  '

# Row 24
Set-TextValue "D24" '157.25'
$ws.Range("E24").Value = '  +1.27%  '

# Row 25
$ws.Range("E25").Value = '  +0.59%  '

# Row 26
Set-TextValue "D26" '8.369'
$ws.Range("E26").Value = '  -0.76%  '

# Row 27
Set-TextValue "D27" '17.64'
$ws.Range("E27").Value = '  -0.39%  '

# Row 28
Set-TextValue "D28" '1.464'
$ws.Range("E28").Value = '  -0.96%  '

# Row 29
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D29" '0.05693'
$ws.Range("E29").Value = '  -2.99%  '

# Row 30
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D30" '1.256'
$ws.Range("E30").Value = '  -0.50%  '

# Row 31
Set-TextValue "D31" '4.124'
$ws.Range("E31").Value = '  +0.69%  '

# Row 32
Set-TextValue "D32" '4.012'
$ws.Range("E32").Value = '  -0.65%  '

# Row 33
Set-TextValue "D33" '1.845'
$ws.Range("E33").Value = '  -2.75%  '

# Row 34
$ws.Range("E34").Value = '  -1.19%  '

# Row 35
Set-TextValue "D35" '0.7144'
$ws.Range("E35").Value = '  -1.08%  '

# Row 36
Set-TextValue "D36" '2.591'
$ws.Range("E36").Value = '  +0.14%  '

# Row 37
$ws.Range("D37").Value = '1.256.25'
$ws.Range("E37").Value = '  +1.51%  '

# Row 38
$ws.Range("E38").Value = '  +1.66%  '

# Row 39
Set-TextValue "D39" '2.784'
$ws.Range("E39").Value = '  -0.35%  '

# Row 40
Set-TextValue "D40" '0.9132'
$ws.Range("E40").Value = '  +0.74%  '

# Row 41
Set-TextValue "D41" '6.209'
$ws.Range("E41").Value = '  +1.09%  '

# Row 42
$ws.Range("E42").Value = '  +0.04%  '

# Row 43
$ws.Range("D43").Value = '2.008.51'
$ws.Range("E43").Value = '  -6.94%  '

# Row 44
Set-TextValue "D44" '101.27'
$ws.Range("E44").Value = '  -0.62%  '

# Row 45
Set-TextValue "D45" '66.40'
$ws.Range("E45").Value = '  -0.87%  '

# Row 46
$ws.Range("E46").Value = '  -3.34%  '

# Row 47
$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D47" '0.4021'
$ws.Range("E47").Value = '  -0.25%  '

# Row 48
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D48" '9.117'
$ws.Range("E48").Value = '  -0.40%  '

# Row 49
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue "D49" '0.00000000116'
$ws.Range("E49").Value = '  -1.40%  '

# Row 50
$ws.Range("E50").Value = '  -1.37%  '

# Row 51
$ws.Range("E51").Value = '  +1.15%  '
